$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the marking scheme values (row 11) and the resulting totals (row 12)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "70 / 112"
